# "Reorganización completa": the old "Datos" sheet becomes "pesajes", the
# bold/white-on-blue, centered header styling applied to row 1 is stripped
# back to the plain default style, and the hard-coded 20-char column widths
# are reset to the worksheet's normal default width.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet "Datos" -> "pesajes"
$ws.Name = "pesajes"

# Drop the custom bold/white-on-blue centered header style (fontId 1 /
# fillId 2 / cellXfs[1]) from A1:E1 so the cells fall back to the default
# "Normal" style again.
$ws.Range("A1:E1").ClearFormats()

# Columns A:E previously had a fixed custom width of 20 characters; restore
# them to the sheet's standard/default column width.
$ws.Columns("A:E").ColumnWidth = 8.43
